$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the part number for the ATMEGA328P row (row 10, column C):
# old value "ATMEGA328P-15AZCT-ND" -> new value "ATMEGA328P-AU-ND"
$ws.Range("C10").Value = "ATMEGA328P-AU-ND"

# Update the selected cell in the sheet view from F18 to C14
$ws.Range("C14").Select()
